# Generate Report for Archive
#
# The "b3ec0fd3-...md" file (row 4 in each status sheet) has moved on to
# translation, so its Status column flips from "Ready for handoff" to
# "In Translation" on the Overview sheet (columns B and C) as well as on
# each per-locale sheet (column B).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B4").Value = "In Translation"
$overview.Range("C4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B4").Value = "In Translation"
